$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Papa" (Terminal Hortofrutícola
# Agro Chillán) dated 2021-09-10. It belongs right above the existing
# 2021-08-30 record (row 168), so insert a new row there which pushes all
# the rows below (168-186) down by one (to 169-187).
$ws.Rows.Item(168).Insert()

$ws.Cells.Item(168,1).Value = 7
$ws.Cells.Item(168,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(168,3).Value = "Ñuble"
$ws.Cells.Item(168,4).Value = "2021-09-10"
$ws.Cells.Item(168,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(168,5).Value = 16
$ws.Cells.Item(168,6).Value = 100114001
$ws.Cells.Item(168,7).Value = "Papa"
$ws.Cells.Item(168,8).Value = "Patagonia"
$ws.Cells.Item(168,9).Value = "1a (guarda)"
$ws.Cells.Item(168,10).Value = 300
$ws.Cells.Item(168,11).Value = 7500
$ws.Cells.Item(168,12).Value = 8000
$ws.Cells.Item(168,13).Value = 7750
$ws.Cells.Item(168,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(168,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(168,16).Value = 310
$ws.Cells.Item(168,17).Value = 25
$ws.Cells.Item(168,18).Value = "Hortaliza"
